# Generate Report for Handoff
#
# Updates the "Status" and datetime columns to reflect a new handoff, and
# narrows columns E/F (Overview) and C (zh-cn / de-de) which previously held
# the long "Handed back..." status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Latest handoff / HO xliff generate datetimes
$zhcn.Range("H2").Value = "2016-10-19 17:35:28"
$overview.Range("G2").Value = "2016-10-19 17:35:40"
$dede.Range("H2").Value = "2016-10-19 17:35:40"

# Narrow columns that used to hold the long status string (target stored
# width ~17.216 chars; ColumnWidth snaps to the nearest 1/6-char pixel
# grid, so feed it a value that lands on the closest achievable cell).
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
